$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.893344666666666
$ws.Range("H2").Value = 5.680033999999999
$ws.Range("I2").Value = 0.05525983881677096
$ws.Range("J2").Value = 0.05525983881677096
$ws.Range("Q2").Value = 0.5199641968875555
$ws.Range("R2").Value = 4.679677771988
$ws.Range("S2").Value = 0.05525983881677096
$ws.Range("T2").Value = 0.05525983881677096

# Row 3
$ws.Range("G3").Value = 4.159773333333334
$ws.Range("I3").Value = 0.1214086415227279
$ws.Range("J3").Value = 0.1214086415227279
$ws.Range("Q3").Value = 1.142387457804445
$ws.Range("S3").Value = 0.1214086415227279
$ws.Range("T3").Value = 0.1214086415227279

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.240212333333333
$ws.Range("H4").Value = 3.720637
$ws.Range("I4").Value = 0.03619728348733726
$ws.Range("J4").Value = 0.03619728348733727
$ws.Range("Q4").Value = 0.3405962058704444
$ws.Range("R4").Value = 3.065365852834
$ws.Range("S4").Value = 0.03619728348733726
$ws.Range("T4").Value = 0.03619728348733727

# Row 5
$ws.Range("G5").Value = 26.96925
$ws.Range("H5").Value = 80.90774999999999
$ws.Range("I5").Value = 0.7871342361731639
$ws.Range("J5").Value = 0.7871342361731638
$ws.Range("Q5").Value = 7.4064932095
$ws.Range("R5").Value = 66.65843888549999
$ws.Range("S5").Value = 0.7871342361731639
$ws.Range("T5").Value = 0.7871342361731638
